# Fixed issues with individual counties
# 1. Swap the month/year data in columns A and B (header row stays the same,
#    but the underlying data had month and year reversed).
# 2. Add new "grade_*" and survey-frequency columns (H:P) with header labels
#    and a constant set of values replicated on every data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: swap A/B data values for rows 2-13 (month should be in A, year in B) ---
for ($r = 2; $r -le 13; $r++) {
    $colA = $ws.Cells.Item($r, 1).Value2
    $colB = $ws.Cells.Item($r, 2).Value2
    $ws.Cells.Item($r, 1).Value2 = $colB
    $ws.Cells.Item($r, 2).Value2 = $colA
}

# --- Step 2: add new headers in H1:P1 ---
$headers = @("grade_total", "grade_distance", "grade_visitation", "grade_encounters", "NEVER", "RARELY", "SOMETIMES", "FREQUENTLY", "ALWAYS")
$col = 8
foreach ($h in $headers) {
    $ws.Cells.Item(1, $col).Value2 = $h
    $col++
}

# --- Step 3: fill H2:P13 with the new constant values on every data row ---
$values = @(1, 0, 0, 3, 1.076, 1.067, 1.118, 1.211, 1.528)
for ($r = 2; $r -le 13; $r++) {
    $col = 8
    foreach ($v in $values) {
        $ws.Cells.Item($r, $col).Value2 = $v
        $col++
    }
}
